# Add a "gene_id" column (Ensembl ID) next to the existing "symbol" column
# on every worksheet of the workbook, as described in the commit:
#   "Add files via upload / with ensembl ids added"

$wb = $excel.ActiveWorkbook

# Mapping from gene symbol (column F "symbol") to Ensembl gene id
$geneMap = @{
    "Trim15" = "ENSMUSG00000050747"
    "Satb2" = "ENSMUSG00000038331"
    "Trim14" = "ENSMUSG00000039853"
    "Bcl11b" = "ENSMUSG00000048251"
    "Tcf19" = "ENSMUSG00000050410"
    "Zbtb7b" = "ENSMUSG00000028042"
    "Zfp646" = "ENSMUSG00000049739"
    "Plagl1" = "ENSMUSG00000019817"
    "Ulk1" = "ENSMUSG00000029512"
    "Trim27" = "ENSMUSG00000021326"
    "Ngfr" = "ENSMUSG00000000120"
    "Six1" = "ENSMUSG00000051367"
    "Ahr" = "ENSMUSG00000019256"
    "Mcph1" = "ENSMUSG00000039842"
    "Prdm10" = "ENSMUSG00000042496"
    "Chd6" = "ENSMUSG00000057133"
    "Aco1" = "ENSMUSG00000028405"
    "Ncoa6" = "ENSMUSG00000038369"
    "Pcbp3" = "ENSMUSG00000001120"
    "Hoxa9" = "ENSMUSG00000038227"
    "Tug1" = "ENSMUSG00000056579"
    "Dnm1l" = "ENSMUSG00000022789"
    "Pcbp1" = "ENSMUSG00000051695"
    "Cdk4" = "ENSMUSG00000006728"
    "Il10" = "ENSMUSG00000016529"
    "Aff2" = "ENSMUSG00000031189"
    "Rom1" = "ENSMUSG00000071648"
    "Prdm4" = "ENSMUSG00000035529"
    "Trim68" = "ENSMUSG00000073968"
    "Nkapl" = "ENSMUSG00000059395"
    "Hoxd9" = "ENSMUSG00000043342"
    "Shh" = "ENSMUSG00000002633"
    "Igf2bp2" = "ENSMUSG00000033581"
    "Prdm9" = "ENSMUSG00000051977"
    "Zfp142" = "ENSMUSG00000026135"
    "Prdm12" = "ENSMUSG00000079466"
    "Atp7a" = "ENSMUSG00000033792"
    "Isl1" = "ENSMUSG00000042258"
    "Cdh1" = "ENSMUSG00000000303"
    "Ankra2" = "ENSMUSG00000021661"
    "Kmt2c" = "ENSMUSG00000038056"
    "Lmx1b" = "ENSMUSG00000038765"
    "Cdk6" = "ENSMUSG00000040274"
    "Aplnr" = "ENSMUSG00000044338"
    "Trem2" = "ENSMUSG00000023992"
    "Tgfb1" = "ENSMUSG00000002603"
    "Trim34b" = "ENSMUSG00000090215"
    "Cdkn2a" = "ENSMUSG00000044303"
    "Hrg" = "ENSMUSG00000022877"
    "Trim30a" = "ENSMUSG00000030921"
    "Igf2bp3" = "ENSMUSG00000029814"
    "Myb" = "ENSMUSG00000019982"
    "Alpk2" = "ENSMUSG00000032845"
    "Dnmt3l" = "ENSMUSG00000000730"
    "Trim40" = "ENSMUSG00000073399"
    "Fgf4" = "ENSMUSG00000050917"
    "Nfe2" = "ENSMUSG00000058794"
    "Il2rg" = "ENSMUSG00000031304"
    "Trim38" = "ENSMUSG00000064140"
}

foreach ($ws in $wb.Worksheets) {
    # Determine the last used row by scanning column F (symbol) starting at row 1
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row
    if ($lastRow -lt 1) { $lastRow = 1 }

    # Header for the new column
    $ws.Cells.Item(1, 7).Value2 = "gene_id"

    if ($lastRow -ge 2) {
        for ($r = 2; $r -le $lastRow; $r++) {
            $symbol = $ws.Cells.Item($r, 6).Value2
            if ($symbol -ne $null -and $geneMap.ContainsKey($symbol)) {
                $ws.Cells.Item($r, 7).Value2 = $geneMap[$symbol]
            }
        }
    }
}

$wb.Save()
